$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 4.39364806774175
$ws.Range("C2").Value = 0.226742093337631
$ws.Range("D2").Value = 19.3772933956263
$ws.Range("E2").Value = 0.0000000000000000000000000000000000183608598096176

# Row 3 - depression_mc
$ws.Range("B3").Value = -0.205137240240449
$ws.Range("C3").Value = 0.457046094563869
$ws.Range("D3").Value = -0.448832716612969
$ws.Range("E3").Value = 0.654533259868536

# Row 4 - anhedonia_mc
$ws.Range("B4").Value = 0.194228571524822
$ws.Range("C4").Value = 0.453668695401615
$ws.Range("D4").Value = 0.428128661936612
$ws.Range("E4").Value = 0.669487412448087

# Row 5 - depression_mc:anhedonia_mc
$ws.Range("B5").Value = -2.13965244830674
$ws.Range("C5").Value = 0.912753451076618
$ws.Range("D5").Value = -2.34417349590183
$ws.Range("E5").Value = 0.0210670244572095
